$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new data row (row 17) with the latest "Random" method scan result,
# following the same layout as the existing rows above it.
$row = 17

# Copy the formatting (date/time number format) from the cell directly above
# so the new date cell reuses the existing style instead of creating a new one.
$ws.Range("A16").Copy()
$ws.Range("A$row").PasteSpecial(-4122)  # xlPasteFormats

$ws.Cells.Item($row, 1).Value = 42620.885636574072
$ws.Cells.Item($row, 2).Value = 56
$ws.Cells.Item($row, 3).Value = 0
$ws.Cells.Item($row, 4).Value = 0
$ws.Cells.Item($row, 5).Value = 0
$ws.Cells.Item($row, 6).Value = 0
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Value = 0
$ws.Cells.Item($row, 9).Value = 0
$ws.Cells.Item($row, 10).Value = 0
$ws.Cells.Item($row, 11).Value = 0
$ws.Cells.Item($row, 12).Value = 0
$ws.Cells.Item($row, 13).Value = 0
$ws.Cells.Item($row, 14).Value = "Random"
